$wb = $excel.ActiveWorkbook

# --- i-type: selection becomes "select all cells" (as if the user hit Ctrl+A) ---
$iType = $wb.Worksheets.Item("i-type")
$iType.Select()
$iType.Cells.Select()

# --- s-type: replace the empty placeholder sheet with a populated copy of i-type ---
$sType = $wb.Worksheets.Item("s-type")
$sType.Delete()

# Duplicate i-type right after itself (Excel names the duplicate "i-type (2)")
$iType.Copy($iType)
$newSheet = $wb.Worksheets.Item("i-type (2)")
$newSheet.Name = "s-type"

# Move the freshly named sheet to the end, restoring the original tab order
$newSheet.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# Re-fetch by name (the old reference becomes stale after Move) and set it active
$newSheet = $wb.Worksheets.Item("s-type")
$newSheet.Select()
$newSheet.Range("AH12").Select()
